$d = $word.ActiveDocument

# --- 1. "Game End:" bullet: remove the trailing paragraph spacing (it is no
#        longer the last item in the Additional Rules list) ---------------
$gameEnd = $d.Paragraphs.Item(62)
$gameEnd.Format.SpaceAfter = 0
$gameEnd.Format.SpaceAfterAuto = $false

# --- 2. Replace the trailing empty paragraph with a new "Optional 1 Joker:"
#        bullet that continues the same numbered/bulleted list (numId 5) ---
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:pPr>' +
         '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr>' +
         '<w:spacing w:after="240" w:before="0" w:beforeAutospacing="0" w:lineRule="auto"/>' +
         '<w:ind w:left="720" w:hanging="360"/>' +
         '<w:rPr><w:u w:val="none"/></w:rPr>' +
       '</w:pPr>' +
       '<w:r><w:rPr><w:b w:val="1"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Optional 1 Joker:</w:t></w:r>' +
       '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> If you desire a flair of whimsy, add 1 Joker to each deck. When a Joker is played, choose a face card (excluding Aces) from any suit and resolve its effect.</w:t></w:r>' +
       '</w:p>'
$last.Range.InsertXML($xml)

# Make sure the (otherwise default) before-spacing values are still written
# out explicitly, matching the rest of the list's paragraphs.
$newLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$newLast.Format.SpaceBefore = 0
$newLast.Format.SpaceBeforeAuto = $false

Write-Output "done"
